$d = $word.ActiveDocument

# 1) Bump the application fee from $75 to $100.
$feeRng = $d.Content
$feeRng.Find.Execute("75", $true, $false, $false, $false, $false, `
                      $true, 1, $false, "100", 2) | Out-Null

# 2) Re-affirm the tenancy-history permission sentence (the "any and all"
#    clause) so Word's stale grammar-check markers around it are cleared.
$permitRng = $d.Content
$permitText = "I the undersigned permit the landlord named above or his or " + `
    "her representative(s) to contact any current and/or previous landlord " + `
    "and to discuss with him or her or his or her representatives any and " + `
    "all information pertaining to my tenancy for the purpose of, but not " + `
    "limited to, establishing my rental payment history and care of the " + `
    "rental property."
$permitRng.Find.Execute($permitText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $permitText, 2) | Out-Null

# 3) Re-affirm the discrimination notice sentence (the "on the basis of"
#    clause) so Word's stale grammar-check markers around it are cleared.
$discrimRng = $d.Content
$discrimText = "It is against the law to discriminate against tenants on " + `
    "the basis of race, ethnicity, sex, sexual orientation, creed, " + `
    "national origin/ancestry, age, marital status, student status, " + `
    "disability, or Vietnam-Era veteran status."
$discrimRng.Find.Execute($discrimText, $true, $false, $false, $false, $false, `
                          $true, 1, $false, $discrimText, 2) | Out-Null
